# Updates NATMI TPM-derived ligand-receptor statistics (per commit: 'update scripts wuth new tpm')
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 2.629732666666667
$ws.Cells.Item(2, 8).Value = 7.889198
$ws.Cells.Item(2, 9).Value = 0.07156737804735891
$ws.Cells.Item(2, 10).Value = 0.07156737804735891
$ws.Cells.Item(2, 13).Value = 0.484733
$ws.Cells.Item(2, 14).Value = 1.454199
$ws.Cells.Item(2, 15).Value = 0.00792098608860474
$ws.Cells.Item(2, 16).Value = 0.00792098608860474
$ws.Cells.Item(2, 17).Value = 1.274718204711333
$ws.Cells.Item(2, 18).Value = 11.472463842402
$ws.Cells.Item(2, 19).Value = 0.0005668842059110462
$ws.Cells.Item(2, 20).Value = 0.0005668842059110462

$ws.Cells.Item(3, 7).Value = 2.629732666666667
$ws.Cells.Item(3, 8).Value = 7.889198
$ws.Cells.Item(3, 9).Value = 0.07156737804735891
$ws.Cells.Item(3, 10).Value = 0.07156737804735891
$ws.Cells.Item(3, 15).Value = 0.1147190689515559
$ws.Cells.Item(3, 16).Value = 0.1147190689515559
$ws.Cells.Item(3, 17).Value = 18.46165161563134
$ws.Cells.Item(3, 18).Value = 166.154864540682
$ws.Cells.Item(3, 19).Value = 0.008210142976897034
$ws.Cells.Item(3, 20).Value = 0.008210142976897034

$ws.Cells.Item(4, 7).Value = 2.629732666666667
$ws.Cells.Item(4, 8).Value = 7.889198
$ws.Cells.Item(4, 9).Value = 0.07156737804735891
$ws.Cells.Item(4, 10).Value = 0.07156737804735891
$ws.Cells.Item(4, 13).Value = 53.289524
$ws.Cells.Item(4, 14).Value = 159.868572
$ws.Cells.Item(4, 15).Value = 0.8708001689019901
$ws.Cells.Item(4, 16).Value = 0.8708001689019901
$ws.Cells.Item(4, 17).Value = 140.1372020539173
$ws.Cells.Item(4, 18).Value = 1261.234818485256
$ws.Cells.Item(4, 19).Value = 0.06232088489151272
$ws.Cells.Item(4, 20).Value = 0.06232088489151272

$ws.Cells.Item(5, 7).Value = 2.629732666666667
$ws.Cells.Item(5, 8).Value = 7.889198
$ws.Cells.Item(5, 9).Value = 0.07156737804735891
$ws.Cells.Item(5, 10).Value = 0.07156737804735891
$ws.Cells.Item(5, 13).Value = 0.4014323333333333
$ws.Cells.Item(5, 14).Value = 1.204297
$ws.Cells.Item(5, 15).Value = 0.006559776057849319
$ws.Cells.Item(5, 16).Value = 0.006559776057849319
$ws.Cells.Item(5, 17).Value = 1.055659720422889
$ws.Cells.Item(5, 18).Value = 9.500937483806
$ws.Cells.Item(5, 19).Value = 0.0004694659730381159
$ws.Cells.Item(5, 20).Value = 0.0004694659730381159

$ws.Cells.Item(6, 9).Value = 0.493312042610523
$ws.Cells.Item(6, 10).Value = 0.493312042610523
$ws.Cells.Item(6, 13).Value = 0.484733
$ws.Cells.Item(6, 14).Value = 1.454199
$ws.Cells.Item(6, 15).Value = 0.00792098608860474
$ws.Cells.Item(6, 16).Value = 0.00792098608860474
$ws.Cells.Item(6, 17).Value = 8.786598845396336
$ws.Cells.Item(6, 18).Value = 79.07938960856701
$ws.Cells.Item(6, 19).Value = 0.003907517826859141
$ws.Cells.Item(6, 20).Value = 0.003907517826859141

$ws.Cells.Item(7, 9).Value = 0.493312042610523
$ws.Cells.Item(7, 10).Value = 0.493312042610523
$ws.Cells.Item(7, 15).Value = 0.1147190689515559
$ws.Cells.Item(7, 16).Value = 0.1147190689515559
$ws.Cells.Item(7, 19).Value = 0.05659229823086947
$ws.Cells.Item(7, 20).Value = 0.05659229823086947

$ws.Cells.Item(8, 9).Value = 0.493312042610523
$ws.Cells.Item(8, 10).Value = 0.493312042610523
$ws.Cells.Item(8, 13).Value = 53.289524
$ws.Cells.Item(8, 14).Value = 159.868572
$ws.Cells.Item(8, 15).Value = 0.8708001689019901
$ws.Cells.Item(8, 16).Value = 0.8708001689019901
$ws.Cells.Item(8, 17).Value = 965.9620245580975
$ws.Cells.Item(8, 18).Value = 8693.658221022877
$ws.Cells.Item(8, 19).Value = 0.4295762100266292
$ws.Cells.Item(8, 20).Value = 0.4295762100266292

$ws.Cells.Item(9, 9).Value = 0.493312042610523
$ws.Cells.Item(9, 10).Value = 0.493312042610523
$ws.Cells.Item(9, 13).Value = 0.4014323333333333
$ws.Cells.Item(9, 14).Value = 1.204297
$ws.Cells.Item(9, 15).Value = 0.006559776057849319
$ws.Cells.Item(9, 16).Value = 0.006559776057849319
$ws.Cells.Item(9, 17).Value = 7.276634511311223
$ws.Cells.Item(9, 18).Value = 65.48971060180101
$ws.Cells.Item(9, 19).Value = 0.003236016526165252
$ws.Cells.Item(9, 20).Value = 0.003236016526165252

$ws.Cells.Item(10, 7).Value = 7.550656333333333
$ws.Cells.Item(10, 8).Value = 22.651969
$ws.Cells.Item(10, 9).Value = 0.2054888252189962
$ws.Cells.Item(10, 10).Value = 0.2054888252189962
$ws.Cells.Item(10, 13).Value = 0.484733
$ws.Cells.Item(10, 14).Value = 1.454199
$ws.Cells.Item(10, 15).Value = 0.00792098608860474
$ws.Cells.Item(10, 16).Value = 0.00792098608860474
$ws.Cells.Item(10, 17).Value = 3.660052296425667
$ws.Cells.Item(10, 18).Value = 32.940470667831
$ws.Cells.Item(10, 19).Value = 0.0016276741259234
$ws.Cells.Item(10, 20).Value = 0.0016276741259234

$ws.Cells.Item(11, 7).Value = 7.550656333333333
$ws.Cells.Item(11, 8).Value = 22.651969
$ws.Cells.Item(11, 9).Value = 0.2054888252189962
$ws.Cells.Item(11, 10).Value = 0.2054888252189962
$ws.Cells.Item(11, 15).Value = 0.1147190689515559
$ws.Cells.Item(11, 16).Value = 0.1147190689515559
$ws.Cells.Item(11, 17).Value = 53.00827284168567
$ws.Cells.Item(11, 18).Value = 477.074455575171
$ws.Cells.Item(11, 19).Value = 0.02357348670907224
$ws.Cells.Item(11, 20).Value = 0.02357348670907225

$ws.Cells.Item(12, 7).Value = 7.550656333333333
$ws.Cells.Item(12, 8).Value = 22.651969
$ws.Cells.Item(12, 9).Value = 0.2054888252189962
$ws.Cells.Item(12, 10).Value = 0.2054888252189962
$ws.Cells.Item(12, 13).Value = 53.289524
$ws.Cells.Item(12, 14).Value = 159.868572
$ws.Cells.Item(12, 15).Value = 0.8708001689019901
$ws.Cells.Item(12, 16).Value = 0.8708001689019901
$ws.Cells.Item(12, 17).Value = 402.3708818909187
$ws.Cells.Item(12, 18).Value = 3621.337937018268
$ws.Cells.Item(12, 19).Value = 0.1789397037081734
$ws.Cells.Item(12, 20).Value = 0.1789397037081734

$ws.Cells.Item(13, 7).Value = 7.550656333333333
$ws.Cells.Item(13, 8).Value = 22.651969
$ws.Cells.Item(13, 9).Value = 0.2054888252189962
$ws.Cells.Item(13, 10).Value = 0.2054888252189962
$ws.Cells.Item(13, 13).Value = 0.4014323333333333
$ws.Cells.Item(13, 14).Value = 1.204297
$ws.Cells.Item(13, 15).Value = 0.006559776057849319
$ws.Cells.Item(13, 16).Value = 0.006559776057849319
$ws.Cells.Item(13, 17).Value = 3.031077590088111
$ws.Cells.Item(13, 18).Value = 27.279698310793
$ws.Cells.Item(13, 19).Value = 0.001347960675827154
$ws.Cells.Item(13, 20).Value = 0.001347960675827155

$ws.Cells.Item(14, 7).Value = 8.437784666666667
$ws.Cells.Item(14, 8).Value = 25.313354
$ws.Cells.Item(14, 9).Value = 0.2296317541231219
$ws.Cells.Item(14, 10).Value = 0.2296317541231219
$ws.Cells.Item(14, 13).Value = 0.484733
$ws.Cells.Item(14, 14).Value = 1.454199
$ws.Cells.Item(14, 15).Value = 0.00792098608860474
$ws.Cells.Item(14, 16).Value = 0.00792098608860474
$ws.Cells.Item(14, 17).Value = 4.090072674827334
$ws.Cells.Item(14, 18).Value = 36.810654073446
$ws.Cells.Item(14, 19).Value = 0.001818909929911152
$ws.Cells.Item(14, 20).Value = 0.001818909929911152

$ws.Cells.Item(15, 7).Value = 8.437784666666667
$ws.Cells.Item(15, 8).Value = 25.313354
$ws.Cells.Item(15, 9).Value = 0.2296317541231219
$ws.Cells.Item(15, 10).Value = 0.2296317541231219
$ws.Cells.Item(15, 15).Value = 0.1147190689515559
$ws.Cells.Item(15, 16).Value = 0.1147190689515559
$ws.Cells.Item(15, 17).Value = 59.23622689798734
$ws.Cells.Item(15, 18).Value = 533.126042081886
$ws.Cells.Item(15, 19).Value = 0.02634314103471715
$ws.Cells.Item(15, 20).Value = 0.02634314103471715

$ws.Cells.Item(16, 7).Value = 8.437784666666667
$ws.Cells.Item(16, 8).Value = 25.313354
$ws.Cells.Item(16, 9).Value = 0.2296317541231219
$ws.Cells.Item(16, 10).Value = 0.2296317541231219
$ws.Cells.Item(16, 13).Value = 53.289524
$ws.Cells.Item(16, 14).Value = 159.868572
$ws.Cells.Item(16, 15).Value = 0.8708001689019901
$ws.Cells.Item(16, 16).Value = 0.8708001689019901
$ws.Cells.Item(16, 17).Value = 449.6455285011654
$ws.Cells.Item(16, 18).Value = 4046.809756510488
$ws.Cells.Item(16, 19).Value = 0.1999633702756748
$ws.Cells.Item(16, 20).Value = 0.1999633702756748

$ws.Cells.Item(17, 7).Value = 8.437784666666667
$ws.Cells.Item(17, 8).Value = 25.313354
$ws.Cells.Item(17, 9).Value = 0.2296317541231219
$ws.Cells.Item(17, 10).Value = 0.2296317541231219
$ws.Cells.Item(17, 13).Value = 0.4014323333333333
$ws.Cells.Item(17, 14).Value = 1.204297
$ws.Cells.Item(17, 15).Value = 0.006559776057849319
$ws.Cells.Item(17, 16).Value = 0.006559776057849319
$ws.Cells.Item(17, 17).Value = 3.387199586904222
$ws.Cells.Item(17, 18).Value = 30.484796282138
$ws.Cells.Item(17, 19).Value = 0.001506332882818796
$ws.Cells.Item(17, 20).Value = 0.001506332882818796
